$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark: it currently sits in the "Weeding" section
#    paragraph (an artifact of the previous edit). Remove it there and
#    recreate it at the end of the new "grammar" paragraph we add below
#    (reflecting that the new paragraph is this commit's last edit point).
# ---------------------------------------------------------------------------
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

# ---------------------------------------------------------------------------
# 2) At the end of the "Parsing" section body paragraph ("...generated by
#    the scanner.") append a trailing space run, then start a brand new
#    paragraph with the grammar blurb.
# ---------------------------------------------------------------------------
$parsingBody = $d.Paragraphs.Item(20)
$parsingBody.Range.InsertAfter(" ")

$afterParsing = $d.Paragraphs.Item(20).Range
$afterParsing.Collapse(0)
$afterParsing.InsertParagraphAfter()

$grammarPara = $d.Paragraphs.Item(21).Range
# Append a trailing placeholder character so the bookmark we add next does
# not land on the (buggy) "last character before the paragraph mark"
# boundary; we strip the placeholder right after anchoring the bookmark.
$grammarPara.InsertAfter("The grammar was developed with extensive reference to the online Java documentation as well as the specific details of the Joos language.X")

$grammarRange = $d.Paragraphs.Item(21).Range
$bookmarkPos = $grammarRange.End - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$goBack = $d.Bookmarks.Item("_GoBack")
$placeholder = $d.Range($goBack.End, $goBack.End + 1)
$placeholder.Delete()

# ---------------------------------------------------------------------------
# 3) Append a new "AST" Heading1 section (with a blank paragraph before it)
#    after the final paragraph of the document.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$lastRange.Collapse(0)
$lastRange.InsertParagraphAfter()

$blankRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$blankRange.Collapse(0)
$blankRange.InsertParagraphAfter()

$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$headingPara.Range.InsertAfter("AST")
$headingPara.Style = "Heading1"

$astTexts = @(
  "The compiler will make use of an AST to simplify the parse tree generated by the parser. The AST is still in development, but the current design involves creating a class for each nonterminal in the AST.",
  "These classes will have a meaningful name that should improve code readability for the rest of the parser.",
  "Each class will contain a parent pointer and several child pointers. The class will have a child pointer for every possible type of child that that class could have. If a child is missing from a particular object, the pointer for that child will be initialized to NULL.",
  "There will also be an epsilon class to represent non-terminals that are reduced to null.",
  "The AST classes will eventually be tagged with attributes which can be used by attribute grammars throughout the remainder of the compiler."
)

foreach ($astText in $astTexts) {
    $precedingRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
    $precedingRange.Collapse(0)
    $precedingRange.InsertParagraphAfter()

    $newBodyPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newBodyPara.Style = "Normal"
    $newBodyPara.Range.InsertAfter($astText)
}
